$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.828.54"
$ws.Range("E2").Value = "  -1.83%  "

# Row 3
$ws.Range("D3").Value = "3.939.24"
$ws.Range("E3").Value = "  -2.37%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.94"
$ws.Range("E5").Value = "  +2.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.98"
$ws.Range("E6").Value = "  +0.62%  "

# Row 7
$ws.Range("D7").Value = "3.934.03"
$ws.Range("E7").Value = "  -2.25%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.686"
$ws.Range("E8").Value = "  -4.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.737"
$ws.Range("E10").Value = "  -4.46%  "

# Row 11
$ws.Range("E11").Value = "  -6.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.40"
$ws.Range("E12").Value = "  +13.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000315"
$ws.Range("E13").Value = "  -4.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.57"
$ws.Range("E14").Value = "  -5.28%  "

# Row 15
$ws.Range("D15").Value = "4.567.02"
$ws.Range("E15").Value = "  -2.40%  "

# Row 16
$ws.Range("D16").Value = "3.937.05"
$ws.Range("E16").Value = "  -2.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.53"
$ws.Range("E17").Value = "  -3.49%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.83"
$ws.Range("E18").Value = "  -2.96%  "

# Row 19
$ws.Range("E19").Value = "  -1.55%  "

# Row 20
$ws.Range("E20").Value = "  -4.86%  "

# Row 21
$ws.Range("D21").Value = "70.694.89"
$ws.Range("E21").Value = "  -2.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "422.87"
$ws.Range("E22").Value = "  -4.90%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.59"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "96.85"
$ws.Range("E24").Value = "  -7.89%  "

# Row 25
$ws.Range("E25").Value = "  +4.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.37"
$ws.Range("E26").Value = "  -5.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.33"
$ws.Range("E27").Value = "  -1.68%  "

# Row 28
$ws.Range("E28").Value = "  +14.51%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.62"
$ws.Range("E29").Value = "  -3.86%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.89"
$ws.Range("E30").Value = "  +1.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.28"
$ws.Range("E31").Value = "  -4.11%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.79"
$ws.Range("E32").Value = "  +14.92%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.75"
$ws.Range("E33").Value = "  +19.34%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.131"
$ws.Range("E34").Value = "  -0.12%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.31"
$ws.Range("E35").Value = "  -3.49%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "681.48"
$ws.Range("E36").Value = "  +0.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "65.00"
$ws.Range("E37").Value = "  -4.27%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.439"
$ws.Range("E38").Value = "  +2.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.150"
$ws.Range("E39").Value = "  -1.72%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0816"
$ws.Range("E40").Value = "  -5.25%  "

# Row 41
$ws.Range("E41").Value = "  -4.19%  "

# Row 42
$ws.Range("E42").Value = "  -0.07%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.07%  "

# Row 44
$ws.Range("E44").Value = "  -4.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.18"
$ws.Range("E45").Value = "  -1.22%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.148"
$ws.Range("E46").Value = "  -5.86%  "

# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.94"
$ws.Range("E47").Value = "  +4.28%  "

# Row 48
$ws.Range("E48").Value = "  -1.05%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.34"
$ws.Range("E49").Value = "  -4.91%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.00"
$ws.Range("E50").Value = "  -2.47%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.61"
$ws.Range("E51").Value = "  -1.35%  "
